$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{
        r = 69
        A = 'Transmission correct'
        B = '8l60xdc4'
        C = 'Training phase'
        D = 1
        E = '[''Purple'', ''Orange'', ''Green'']'
        F = '[[''Red'', ''Red''], [''Blue'', ''''], [''Blue'', '''']]'
        G = '[None, None, None]'
        H = '[''5'', ''5'', ''5'']'
        I = '0.15'
    }
    @{
        r = 70
        A = 'Transmission M&M'
        B = 'zkjjaoa6'
        C = 'Training phase'
        D = 1
        E = '[''Purple'', ''Orange'', ''Green'']'
        F = '[[''Blue'', ''''], [''Blue'', ''''], [''Blue'', '''']]'
        G = '[None, None, None]'
        H = '[''5'', ''5'', ''5'']'
        I = '0.15'
    }
    @{
        r = 71
        A = 'Transmission correct'
        B = 'ygg4hk51'
        C = 'Training phase'
        D = 1
        E = '[''Purple'', ''Orange'', ''Green'']'
        F = '[[''Red'', ''''], [''Red'', ''''], [''Red'', '''']]'
        G = '[None, None, None]'
        H = '[''2'', ''2'', ''2'']'
        I = '0.06'
    }
    @{
        r = 72
        A = 'Transmission M&M'
        B = 'bh4i5uou'
        C = 'Training phase'
        D = 2
        E = '[''Green'', ''Green'', ''Orange'']'
        F = '[[''Red'', ''''], [''Blue'', ''''], [''Yellow'', '''']]'
        G = '[None, None, None]'
        H = '[''2'', ''5'', ''3'']'
        I = '0.10'
    }
    @{
        r = 73
        A = 'Transmission M&M'
        B = '94uigg6z'
        C = 'Training phase'
        D = 1
        E = '[''Purple'', ''Orange'', ''Green'']'
        F = '[[''Red'', ''''], [''Red'', ''''], [''Red'', '''']]'
        G = '[None, None, None]'
        H = '[''2'', ''2'', ''2'']'
        I = '0.06'
    }
    @{
        r = 74
        A = 'Transmission correct'
        B = '8qiqijij'
        C = 'Training phase'
        D = 1
        E = '[''Purple'', ''Orange'', ''Green'']'
        F = '[[''Red'', ''''], [''Red'', ''''], [''Red'', '''']]'
        G = '[None, None, None]'
        H = '[''2'', ''2'', ''2'']'
        I = $null
    }
    @{
        r = 75
        A = 'Transmission correct'
        B = '8qiqijij'
        C = 'Training phase'
        D = 2
        E = '[''Green'', ''Green'', ''Orange'']'
        F = '[[''Red'', ''''], [''Red'', ''''], [''Red'', '''']]'
        G = '[None, None, None]'
        H = '[''2'', ''2'', ''2'']'
        I = '0.12'
    }
    @{
        r = 76
        A = 'Transmission M&M'
        B = '7tfvji4f'
        C = 'Training phase'
        D = 1
        E = '[''Purple'', ''Orange'', ''Green'']'
        F = '[[''Red'', ''''], [''Blue'', ''''], [''Yellow'', '''']]'
        G = '[None, None, None]'
        H = '[''2'', ''5'', ''3'']'
        I = '0.10'
    }
    @{
        r = 77
        A = 'Anomaly no noise'
        B = 'vf3wsvvg'
        C = 'Training phase'
        D = 1
        E = '[''Purple'', ''Orange'', ''Green'']'
        F = '[[''Red'', ''''], [''Blue'', ''''], [''Yellow'', '''']]'
        G = '[None, None, None]'
        H = '[''2'', ''5'', ''3'']'
        I = '0.10'
    }
    @{
        r = 78
        A = 'Anomaly no noise'
        B = 'x7c4pnjr'
        C = 'Training phase'
        D = 1
        E = '[''Purple'', ''Orange'', ''Green'']'
        F = '[[''Red'', ''Blue''], [''Red'', ''''], [''Blue'', ''Blue'']]'
        G = '[None, None, None]'
        H = '[''8'', ''2'', ''10'']'
        I = $null
    }
    @{
        r = 79
        A = 'Anomaly no noise'
        B = 'x7c4pnjr'
        C = 'Training phase'
        D = 2
        E = '[''Green'', ''Green'', ''Orange'']'
        F = '[[''Red'', ''''], [''Blue'', ''''], [''Yellow'', '''']]'
        G = '[None, None, None]'
        H = '[''2'', ''5'', ''3'']'
        I = '0.30'
    }
    @{
        r = 80
        A = 'Transmission correct'
        B = 'xcowwowp'
        C = 'Training phase'
        D = 1
        E = '[''Purple'', ''Orange'', ''Green'']'
        F = '[[''Red'', ''Blue''], [''Blue'', ''Blue''], [''Red'', '''']]'
        G = '[None, None, None]'
        H = '[''8'', ''10'', ''2'']'
        I = '0.20'
    }
    @{
        r = 81
        A = 'Transmission M&M'
        B = 'lssyg1wd'
        C = 'Training phase'
        D = 1
        E = '[''Purple'', ''Orange'', ''Green'']'
        F = '[[''Red'', ''Blue''], [''Blue'', ''Blue''], [''Blue'', '''']]'
        G = '[None, None, None]'
        H = '[''8'', ''10'', ''5'']'
        I = '0.23'
    }
    @{
        r = 82
        A = 'M&M no noise'
        B = '8khfl1x3'
        C = 'Training phase'
        D = 1
        E = '[''Purple'', ''Orange'', ''Green'']'
        F = '[[''Red'', ''Blue''], [''Blue'', ''Blue''], [''Red'', '''']]'
        G = '[None, None, None]'
        H = '[''10'', ''5'', ''3'']'
        I = '0.18'
    }
    @{
        r = 83
        A = 'Anomaly noisy'
        B = 'va0caf9s'
        C = 'Training phase'
        D = 1
        E = '[''Purple'', ''Orange'', ''Green'']'
        F = '[[''Red'', ''Blue''], [''Blue'', ''Blue''], [''Red'', '''']]'
        G = '[None, None, {''index'': 2, ''type'': ''increase'', ''amount'': 0.2, ''before'': 0.25, ''after'': 0.45}]'
        H = '[''8'', ''10'', ''4'']'
        I = '0.22'
    }
    @{
        r = 84
        A = 'M&M noisy'
        B = 'sm0rfajc'
        C = 'Training phase'
        D = 1
        E = '[''Purple'', ''Orange'', ''Green'']'
        F = '[[''Red'', ''Blue''], [''Blue'', ''Blue''], [''Red'', '''']]'
        G = '[None, {''index'': 1, ''type'': ''decrease'', ''amount'': -0.2, ''before'': 0.5, ''after'': 0.3}, None]'
        H = '[''10'', ''3'', ''3'']'
        I = '0.16'
    }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    if ($row.I -ne $null) {
        $icell = $ws.Cells.Item($r, 9)
        $icell.NumberFormat = '@'
        $icell.Value = $row.I
        $icell.ClearFormats()
    }
}

Write-Host "Rows added."
